# Binary Search 2: Square root of Integer - updated index
#
# Adds the GitHub source-code hyperlink for the "Square root of Integer"
# problem to the "Binary Search 2" worksheet, in cell F3 - mirroring the
# pattern already used for the GitHub links in column F of "Binary Search 1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")
$cell = $ws.Range("F3")

$githubUrl = "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/1_square_root_of_integer.js"
$friendlyText = "dsa/1_square_root_of_integer.js at main " + [char]0x00B7 + " ankurnecessary/dsa " + [char]0x00B7 + " GitHub"

# Hyperlinks.Add writes both the cell text (TextToDisplay) and the
# hyperlink's stored display/address; set TextToDisplay to the target URL
# here so the "display" attribute matches the raw GitHub link, then
# overwrite the cell's visible text with the friendly GitHub title - same
# two-step shape Excel itself produces for this kind of link.
$link = $ws.Hyperlinks.Add($cell, $githubUrl, [Type]::Missing, [Type]::Missing, $githubUrl)
$cell.Value2 = $friendlyText

# Match the existing "Hyperlink" cell style (wrap text, no explicit
# horizontal/vertical override) instead of the bespoke style Excel's
# Hyperlinks.Add would otherwise synthesize.
$cell.HorizontalAlignment = 1
$cell.VerticalAlignment = -4107
$cell.WrapText = $true
